$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# --- New backlog row (LALA-027) ---
$ws.Range("A28").Value = "LALA-027"
$ws.Range("B28").Value = "P1"
$ws.Range("C28").Value = "UI  "
$ws.Range("D28").Value = "Tire change time changing total fuel unecessarily"
$ws.Range("E28").Value = "Changing tire time that is still less than refuel time is adjusting fuel when it shouldn't because the time is already accounted by the refuel time."
$ws.Range("F28").Value = "When tire change time is less than refuel time, the total fuel needed should not change"
$ws.Range("G28").Value = "Backlog"
$ws.Range("H28").Value = "Andy"
$ws.Range("I28").Value = "any"

# Row 28 renders as a two-line (30pt) row, same as the other wrapped entries.
$ws.Rows.Item(28).RowHeight = 30

# --- Column width tweaks (Description / Acceptance Criteria got wider) ---
$ws.Columns.Item(5).ColumnWidth = 69
$ws.Columns.Item(6).ColumnWidth = 75.5

# --- Freeze the header row and move the live selection to the new row ---
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("J28").Select()
